$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (J1:M1) - Cards and Goals stats
$ws.Range("J1").Value = "Cards For Argentina"
$ws.Range("K1").Value = "Cards For Opponent"
$ws.Range("L1").Value = "Goal For Argentina "
$ws.Range("M1").Value = "Goal Against Argentina"

# Row 2: Saudi Arabia
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2

# Row 3: Mexico
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 0

# Row 4: Poland
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = 0

# Row 5: Australia
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 1

# Row 6: Netherlands
$ws.Range("J6").Value = 10
$ws.Range("K6").Value = 10
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 5

# Row 7: Croatia
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = 0

# Row 8: France
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 7
$ws.Range("M8").Value = 5

# Match the original file's best-fit column widths on the new columns
$ws.Columns.Item(10).ColumnWidth = 16.166666666666668
$ws.Columns.Item(11).ColumnWidth = 19.385416666666668
$ws.Columns.Item(12).ColumnWidth = 15.830729166666666

# Update selection to mirror the saved view state
$ws.Range("O11").Select() | Out-Null
